$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(51, 8).Value = 18566.834
$ws.Cells.Item(51, 9).Value = 1800
$ws.Cells.Item(51, 10).Value = 20091.092
$ws.Cells.Item(51, 11).Value = 1800
$ws.Cells.Item(51, 12).Value = 20091.092
$ws.Cells.Item(51, 13).Value = -1316
$ws.Cells.Item(51, 14).Value = -21059.092

$ws.Cells.Item(116, 8).Value = 19850430
$ws.Cells.Item(116, 9).Value = 14289116
$ws.Cells.Item(116, 10).Value = 23822796
$ws.Cells.Item(116, 11).Value = 14289116
$ws.Cells.Item(116, 12).Value = 23822796
$ws.Cells.Item(116, 13).Value = -14285674
$ws.Cells.Item(116, 14).Value = -23829680

$ws.Cells.Item(128, 8).Value = 15292.857
$ws.Cells.Item(128, 10).Value = 15292.857
$ws.Cells.Item(128, 12).Value = 15292.857
$ws.Cells.Item(128, 14).Value = -25252.857

$ws.Cells.Item(138, 8).Value = 3423.4504
$ws.Cells.Item(138, 9).Value = 2553.3225
$ws.Cells.Item(138, 10).Value = 3873.0166
$ws.Cells.Item(138, 11).Value = 7659.967500000001
$ws.Cells.Item(138, 12).Value = 11619.0498
$ws.Cells.Item(138, 13).Value = -2519.967500000001
$ws.Cells.Item(138, 14).Value = -21899.0498

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4414.1
$ws.Cells.Item(32, 9).Value = 4307.172
$ws.Cells.Item(32, 10).Value = 15000
$ws.Cells.Item(32, 11).Value = 4307.172
$ws.Cells.Item(32, 12).Value = 15000
$ws.Cells.Item(32, 13).Value = -4020.172
$ws.Cells.Item(32, 14).Value = -15574

$ws.Cells.Item(74, 8).Value = 4446265.5
$ws.Cells.Item(74, 9).Value = 963.28125
$ws.Cells.Item(74, 10).Value = 15388547
$ws.Cells.Item(74, 11).Value = 963.28125
$ws.Cells.Item(74, 12).Value = 15388547
$ws.Cells.Item(74, 13).Value = -89.28125
$ws.Cells.Item(74, 14).Value = -15390295

$ws.Cells.Item(77, 8).Value = 4446265.5
$ws.Cells.Item(77, 9).Value = 963.28125
$ws.Cells.Item(77, 10).Value = 15388547
$ws.Cells.Item(77, 11).Value = 4816.40625
$ws.Cells.Item(77, 12).Value = 76942735
$ws.Cells.Item(77, 13).Value = -448.40625
$ws.Cells.Item(77, 14).Value = -76951471

$ws.Cells.Item(132, 8).Value = 24932694
$ws.Cells.Item(132, 9).Value = 30055272
$ws.Cells.Item(132, 10).Value = 10418726
$ws.Cells.Item(132, 11).Value = 90165816
$ws.Cells.Item(132, 12).Value = 31256178
$ws.Cells.Item(132, 13).Value = -90163286
$ws.Cells.Item(132, 14).Value = -31261238

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 10821632
$ws.Cells.Item(134, 9).Value = 11793655
$ws.Cells.Item(134, 10).Value = 129375
$ws.Cells.Item(134, 11).Value = 35380965
$ws.Cells.Item(134, 12).Value = 388125
$ws.Cells.Item(134, 13).Value = -35378430
$ws.Cells.Item(134, 14).Value = -393195

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 795.17645
$ws.Cells.Item(22, 9).Value = 174.63637
$ws.Cells.Item(22, 10).Value = 1932.8334
$ws.Cells.Item(22, 11).Value = 174.63637
$ws.Cells.Item(22, 12).Value = 1932.8334
$ws.Cells.Item(22, 13).Value = 175.36363
$ws.Cells.Item(22, 14).Value = -2632.8334

$ws.Cells.Item(31, 8).Value = 2318043
$ws.Cells.Item(31, 9).Value = 3206139.8
$ws.Cells.Item(31, 10).Value = 8991.4
$ws.Cells.Item(31, 11).Value = 3206139.8
$ws.Cells.Item(31, 12).Value = 8991.4
$ws.Cells.Item(31, 13).Value = -3205844.8
$ws.Cells.Item(31, 14).Value = -9581.4

$ws.Cells.Item(34, 8).Value = 2318043
$ws.Cells.Item(34, 9).Value = 3206139.8
$ws.Cells.Item(34, 10).Value = 8991.4
$ws.Cells.Item(34, 11).Value = 3206139.8
$ws.Cells.Item(34, 12).Value = 8991.4
$ws.Cells.Item(34, 13).Value = -3205937.8
$ws.Cells.Item(34, 14).Value = -9395.4

$ws.Cells.Item(50, 8).Value = 13027
$ws.Cells.Item(50, 10).Value = 13027
$ws.Cells.Item(50, 12).Value = 13027
$ws.Cells.Item(50, 14).Value = -14277

$ws.Cells.Item(51, 8).Value = 32787.625
$ws.Cells.Item(51, 10).Value = 10383.5
$ws.Cells.Item(51, 12).Value = 10383.5
$ws.Cells.Item(51, 14).Value = -11855.5

$ws.Cells.Item(58, 8).Value = 1519664.1
$ws.Cells.Item(58, 9).Value = 7478.3335
$ws.Cells.Item(58, 10).Value = 3031849.8
$ws.Cells.Item(58, 11).Value = 7478.3335
$ws.Cells.Item(58, 12).Value = 3031849.8
$ws.Cells.Item(58, 13).Value = -7275.3335
$ws.Cells.Item(58, 14).Value = -3032255.8

$ws.Cells.Item(59, 8).Value = 17746.75
$ws.Cells.Item(59, 10).Value = 17746.75
$ws.Cells.Item(59, 12).Value = 17746.75
$ws.Cells.Item(59, 14).Value = -20036.75

$ws.Cells.Item(60, 8).Value = 22512.625
$ws.Cells.Item(60, 10).Value = 13367
$ws.Cells.Item(60, 12).Value = 13367
$ws.Cells.Item(60, 14).Value = -14389

$ws.Cells.Item(61, 8).Value = 32787.625
$ws.Cells.Item(61, 10).Value = 10383.5
$ws.Cells.Item(61, 12).Value = 10383.5
$ws.Cells.Item(61, 14).Value = -11079.5

$ws.Cells.Item(74, 8).Value = 17790.416
$ws.Cells.Item(74, 10).Value = 19290.908
$ws.Cells.Item(74, 12).Value = 19290.908
$ws.Cells.Item(74, 14).Value = -21038.908

$ws.Cells.Item(77, 8).Value = 17790.416
$ws.Cells.Item(77, 10).Value = 19290.908
$ws.Cells.Item(77, 12).Value = 57872.724
$ws.Cells.Item(77, 14).Value = -66608.724

$ws.Cells.Item(99, 8).Value = 9264.1
$ws.Cells.Item(99, 9).Value = 15781.8
$ws.Cells.Item(99, 10).Value = 7960.56
$ws.Cells.Item(99, 11).Value = 15781.8
$ws.Cells.Item(99, 12).Value = 7960.56
$ws.Cells.Item(99, 13).Value = -14283.8
$ws.Cells.Item(99, 14).Value = -10956.56

$ws.Cells.Item(126, 8).Value = 9264.1
$ws.Cells.Item(126, 9).Value = 15781.8
$ws.Cells.Item(126, 10).Value = 7960.56
$ws.Cells.Item(126, 11).Value = 47345.39999999999
$ws.Cells.Item(126, 12).Value = 23881.68
$ws.Cells.Item(126, 13).Value = -44875.39999999999
$ws.Cells.Item(126, 14).Value = -28821.68

$ws.Cells.Item(132, 8).Value = 1518.1
$ws.Cells.Item(132, 9).Value = 1060.2778
$ws.Cells.Item(132, 10).Value = 2695.3572
$ws.Cells.Item(132, 11).Value = 3180.8334
$ws.Cells.Item(132, 12).Value = 8086.071599999999
$ws.Cells.Item(132, 13).Value = -650.8334000000004
$ws.Cells.Item(132, 14).Value = -13146.0716

$ws.Cells.Item(136, 8).Value = 1519664.1
$ws.Cells.Item(136, 9).Value = 7478.3335
$ws.Cells.Item(136, 10).Value = 3031849.8
$ws.Cells.Item(136, 11).Value = 22435.0005
$ws.Cells.Item(136, 12).Value = 9095549.399999999
$ws.Cells.Item(136, 13).Value = -19885.0005
$ws.Cells.Item(136, 14).Value = -9100649.399999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 3487037.2
$ws.Cells.Item(5, 9).Value = 5495078
$ws.Cells.Item(5, 10).Value = 2405784.5
$ws.Cells.Item(5, 11).Value = 16485234
$ws.Cells.Item(5, 12).Value = 7217353.5
$ws.Cells.Item(5, 13).Value = -16485122
$ws.Cells.Item(5, 14).Value = -7217577.5

$ws.Cells.Item(56, 8).Value = 2740
$ws.Cells.Item(56, 9).Value = 2740
$ws.Cells.Item(56, 11).Value = 2740
$ws.Cells.Item(56, 13).Value = -2210

$ws.Cells.Item(134, 8).Value = 3033.3333
$ws.Cells.Item(134, 9).Value = 2550
$ws.Cells.Item(134, 11).Value = 7650
$ws.Cells.Item(134, 13).Value = -2580

$ws.Cells.Item(135, 8).Value = 3487037.2
$ws.Cells.Item(135, 9).Value = 5495078
$ws.Cells.Item(135, 10).Value = 2405784.5
$ws.Cells.Item(135, 11).Value = 49455702
$ws.Cells.Item(135, 12).Value = 21652060.5
$ws.Cells.Item(135, 13).Value = -49453167
$ws.Cells.Item(135, 14).Value = -21657130.5

$ws.Cells.Item(139, 8).Value = 22021.5
$ws.Cells.Item(139, 9).Value = 27293.947
$ws.Cells.Item(139, 10).Value = 1986.2
$ws.Cells.Item(139, 11).Value = 81881.841
$ws.Cells.Item(139, 12).Value = 5958.6
$ws.Cells.Item(139, 13).Value = -76741.841
$ws.Cells.Item(139, 14).Value = -16238.6

$ws.Cells.Item(140, 8).Value = 2555.6765
$ws.Cells.Item(140, 9).Value = 2295
$ws.Cells.Item(140, 11).Value = 6885
$ws.Cells.Item(140, 13).Value = -1705

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 5054036
$ws.Cells.Item(132, 9).Value = 5503691.5
$ws.Cells.Item(132, 10).Value = 4134286
$ws.Cells.Item(132, 11).Value = 16511074.5
$ws.Cells.Item(132, 12).Value = 12402858
$ws.Cells.Item(132, 13).Value = -16508544.5
$ws.Cells.Item(132, 14).Value = -12407918

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1731.5555
$ws.Cells.Item(7, 9).Value = 1614
$ws.Cells.Item(7, 10).Value = 1966.6666
$ws.Cells.Item(7, 11).Value = 1614
$ws.Cells.Item(7, 12).Value = 1966.6666
$ws.Cells.Item(7, 13).Value = -1502
$ws.Cells.Item(7, 14).Value = -2190.6666

$ws.Cells.Item(126, 8).Value = 1731.5555
$ws.Cells.Item(126, 9).Value = 1614
$ws.Cells.Item(126, 10).Value = 1966.6666
$ws.Cells.Item(126, 11).Value = 4842
$ws.Cells.Item(126, 12).Value = 5899.9998
$ws.Cells.Item(126, 13).Value = -2372
$ws.Cells.Item(126, 14).Value = -10839.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 132.21428
$ws.Cells.Item(113, 9).Value = 95.57143
$ws.Cells.Item(113, 10).Value = 168.85715
$ws.Cells.Item(113, 11).Value = 286.71429
$ws.Cells.Item(113, 12).Value = 506.57145
$ws.Cells.Item(113, 13).Value = 1883.28571
$ws.Cells.Item(113, 14).Value = -4846.57145

$ws.Cells.Item(132, 8).Value = 712962.5
$ws.Cells.Item(132, 9).Value = 2091.6482
$ws.Cells.Item(132, 10).Value = 3665810.5
$ws.Cells.Item(132, 11).Value = 6274.944600000001
$ws.Cells.Item(132, 12).Value = 10997431.5
$ws.Cells.Item(132, 13).Value = -3744.944600000001
$ws.Cells.Item(132, 14).Value = -11002491.5

$ws.Cells.Item(136, 8).Value = 853.5
$ws.Cells.Item(136, 9).Value = 713.3019
$ws.Cells.Item(136, 10).Value = 1915
$ws.Cells.Item(136, 11).Value = 2139.9057
$ws.Cells.Item(136, 12).Value = 5745
$ws.Cells.Item(136, 13).Value = 410.0942999999997
$ws.Cells.Item(136, 14).Value = -10845
